$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 is being replaced: the old "Lab 05" entry is removed and a new
# "MLR: Model Comparison" class day is added in its place.

# Copy the cell formatting from row 28 (a similarly-structured data row)
# into the cells of row 29 that need a different style than the old
# "Lab 05" row had.
$ws.Range("D28").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Range("F28").Copy()
$ws.Range("F29").PasteSpecial(-4122)

$ws.Range("G28").Copy()
$ws.Range("G29").PasteSpecial(-4122)

$ws.Range("H28").Copy()
$ws.Range("H29").PasteSpecial(-4122)

$ws.Range("J28").Copy()
$ws.Range("I29:K29").PasteSpecial(-4122)

$ws.Range("F28").Copy()
$ws.Range("L29").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Now set the new values for row 29 (order matches how the strings were
# entered so that new shared-string entries are appended in the same order).
$ws.Range("A29").Value = "x"
$ws.Range("B29").Value = ""
$ws.Range("C29").Value = "W"
$ws.Range("D29").Value = 45588
$ws.Range("E29").Value = "MLR: Model Comparison"
$ws.Range("G29").Value = "/slides/20-comparison.html"
$ws.Range("H29").Value = "/ae/ae-14-comparison.html"
$ws.Range("F29").Value = "/prepare/mlr-comparison-prep.html"
$ws.Range("I29").Value = ""
$ws.Range("J29").Value = ""
$ws.Range("K29").Value = ""
$ws.Range("L29").Value = ""

$ws.Range("A29:L29").RowHeight = 66
